$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 571.4167
$ws.Range("I80").Value = 672.3333
$ws.Range("J80").Value = 268.66666
$ws.Range("K80").Value = 2016.9999
$ws.Range("L80").Value = 805.9999799999999
$ws.Range("M80").Value = -1018.9999
$ws.Range("N80").Value = -2801.99998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 571.4167
$ws.Range("I83").Value = 672.3333
$ws.Range("J83").Value = 268.66666
$ws.Range("K83").Value = 6050.9997
$ws.Range("L83").Value = 2417.99994
$ws.Range("M83").Value = -1058.9997
$ws.Range("N83").Value = -12401.99994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 109991.664
$ws.Range("J87").Value = 109991.664
$ws.Range("L87").Value = 109991.664
$ws.Range("N87").Value = -112487.664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 109991.664
$ws.Range("J90").Value = 109991.664
$ws.Range("L90").Value = 329974.992
$ws.Range("N90").Value = -342454.992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 2567.6667
$ws.Range("J103").Value = 2127.8
$ws.Range("L103").Value = 6383.400000000001
$ws.Range("N103").Value = -7555.400000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5237.303
$ws.Range("I132").Value = 5339.069
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 16017.207
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -13487.207
$ws.Range("N132").Value = -18558.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 69998.91
$ws.Range("J134").Value = 69998.91
$ws.Range("L134").Value = 69998.91
$ws.Range("N134").Value = -80138.91

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5129.64
$ws.Range("I32").Value = 5129.64
$ws.Range("K32").Value = 5129.64
$ws.Range("M32").Value = -4842.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 45070.6
$ws.Range("I45").Value = 71527.336
$ws.Range("K45").Value = 71527.336
$ws.Range("M45").Value = -71150.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3351.72
$ws.Range("I61").Value = 2671.5334
$ws.Range("K61").Value = 2671.5334
$ws.Range("M61").Value = -2459.5334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3859.7058
$ws.Range("I132").Value = 3401.875
$ws.Range("K132").Value = 10205.625
$ws.Range("M132").Value = -7675.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3351.72
$ws.Range("I136").Value = 2671.5334
$ws.Range("K136").Value = 8014.600199999999
$ws.Range("M136").Value = -5464.600199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 21747846
$ws.Range("I20").Value = 26325280
$ws.Range("K20").Value = 26325280
$ws.Range("M20").Value = -26325033

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2237.25
$ws.Range("I105").Value = 1499.5
$ws.Range("J105").Value = 2975
$ws.Range("K105").Value = 1499.5
$ws.Range("L105").Value = 2975
$ws.Range("M105").Value = 247.5
$ws.Range("N105").Value = -6469

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2717.0527
$ws.Range("I134").Value = 2476.5
$ws.Range("K134").Value = 7429.5
$ws.Range("M134").Value = -4894.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2959.25
$ws.Range("I3").Value = 2245.3076
$ws.Range("K3").Value = 6735.9228
$ws.Range("M3").Value = -6623.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 9011.666999999999
$ws.Range("J62").Value = 9011.666999999999
$ws.Range("L62").Value = 27035.001
$ws.Range("N62").Value = -28407.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3571.1428
$ws.Range("J64").Value = 3999.6667
$ws.Range("L64").Value = 11999.0001
$ws.Range("N64").Value = -12539.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 9011.666999999999
$ws.Range("J65").Value = 9011.666999999999
$ws.Range("L65").Value = 81105.003
$ws.Range("N65").Value = -87969.003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 3571.1428
$ws.Range("J67").Value = 3999.6667
$ws.Range("L67").Value = 11999.0001
$ws.Range("N67").Value = -13871.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4707.0835
$ws.Range("J81").Value = 5148.5
$ws.Range("L81").Value = 15445.5
$ws.Range("N81").Value = -17691.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 4707.0835
$ws.Range("J84").Value = 5148.5
$ws.Range("L84").Value = 46336.5
$ws.Range("N84").Value = -57568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 1345.9
$ws.Range("I138").Value = 1154.4445
$ws.Range("J138").Value = 3069
$ws.Range("K138").Value = 3463.3335
$ws.Range("L138").Value = 9207
$ws.Range("M138").Value = 1676.6665
$ws.Range("N138").Value = -19487

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2522.9678
$ws.Range("I139").Value = 1628
$ws.Range("J139").Value = 2834.261
$ws.Range("K139").Value = 4884
$ws.Range("L139").Value = 8502.782999999999
$ws.Range("M139").Value = 256
$ws.Range("N139").Value = -18782.783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100002830
$ws.Range("I80").Value = 142858690
$ws.Range("J80").Value = 5833
$ws.Range("K80").Value = 142858690
$ws.Range("L80").Value = 5833
$ws.Range("M80").Value = -142857692
$ws.Range("N80").Value = -7829

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 100002830
$ws.Range("I83").Value = 142858690
$ws.Range("J83").Value = 5833
$ws.Range("K83").Value = 714293450
$ws.Range("L83").Value = 29165
$ws.Range("M83").Value = -714288458
$ws.Range("N83").Value = -39149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1171.7646
$ws.Range("I102").Value = 1099.1818
$ws.Range("J102").Value = 1304.8334
$ws.Range("K102").Value = 1099.1818
$ws.Range("L102").Value = 1304.8334
$ws.Range("M102").Value = 522.8181999999999
$ws.Range("N102").Value = -4548.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2990.2144
$ws.Range("I113").Value = 2724.6365
$ws.Range("K113").Value = 2724.6365
$ws.Range("M113").Value = -554.6365000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3640.0667
$ws.Range("I126").Value = 2200.5386
$ws.Range("K126").Value = 6601.6158
$ws.Range("M126").Value = -4131.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 68197.8
$ws.Range("J134").Value = 68197.8
$ws.Range("L134").Value = 204593.4
$ws.Range("N134").Value = -209663.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5505.2104
$ws.Range("I7").Value = 4306.7334
$ws.Range("K7").Value = 4306.7334
$ws.Range("M7").Value = -4194.7334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3427.4443
$ws.Range("J46").Value = 4666.3335
$ws.Range("L46").Value = 4666.3335
$ws.Range("N46").Value = -5042.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2902
$ws.Range("I68").Value = 3250.25
$ws.Range("J68").Value = 2437.6667
$ws.Range("K68").Value = 3250.25
$ws.Range("L68").Value = 2437.6667
$ws.Range("M68").Value = -2501.25
$ws.Range("N68").Value = -3935.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2902
$ws.Range("I71").Value = 3250.25
$ws.Range("J71").Value = 2437.6667
$ws.Range("K71").Value = 16251.25
$ws.Range("L71").Value = 12188.3335
$ws.Range("M71").Value = -12507.25
$ws.Range("N71").Value = -19676.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5505.2104
$ws.Range("I126").Value = 4306.7334
$ws.Range("K126").Value = 12920.2002
$ws.Range("M126").Value = -10450.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 81706.5
$ws.Range("J46").Value = 81706.5
$ws.Range("L46").Value = 81706.5
$ws.Range("N46").Value = -82168.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 30309400
$ws.Range("I132").Value = 41672300
$ws.Range("J132").Value = 8332.666999999999
$ws.Range("K132").Value = 125016900
$ws.Range("L132").Value = 24998.001
$ws.Range("M132").Value = -125014370
$ws.Range("N132").Value = -30058.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 88498.25
$ws.Range("J133").Value = 88498.25
$ws.Range("L133").Value = 88498.25
$ws.Range("N133").Value = -98618.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 81706.5
$ws.Range("J134").Value = 81706.5
$ws.Range("L134").Value = 245119.5
$ws.Range("N134").Value = -250189.5
